$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.288.36"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "2.434.17"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'563.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "'142.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "2.430.05"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -5.36%  "
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "'26.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("E15").Value = "  -6.11%  "
$ws.Range("D16").Value = "2.877.12"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "62.151.18"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "2.434.27"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("E19").Value = "  -4.34%  "
$ws.Range("D20").Value = "'7.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").Value = "'324.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  +3.52%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'64.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").Value = "'619.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "'8.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "2.549.70"
$ws.Range("D29").Value = "0.0₃0954"
$ws.Range("E29").Value = "  -9.37%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("E32").Value = "  -4.99%  "
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("E34").Value = "  -7.48%  "
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -6.96%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "'18.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "'146.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").Value = "'5.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.41%  "
$ws.Range("E42").Value = "  -6.62%  "
$ws.Range("D43").Value = "'42.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("D45").Value = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.52%  "
$ws.Range("D46").Value = "'145.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.38%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").Value = "'19.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.69%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("E51").Value = "  -5.10%  "
